$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

$shp = $s.Shapes.AddShape(1, 10620892, 3054012, 963827, 815546)
$shp.Name = "Rectangle 52"

$shp.TextFrame.TextRange.Text = "for.`rglobal"
